{"js": "// Paragraph 2: split the single run into 5 runs, underlining\n// \"mailing transactional\" and bolding \"Column enhanced fixed\".\nconst body = context.document.body;\n\nconst underlineHits = body.search(\"mailing transactional\", { matchCase: true });\nconst boldHits = body.search(\"Column enhanced fixed\", { matchCase: true });\nunderlineHits.load(\"items\");\nboldHits.load(\"items\");\nawait context.sync();\n\nunderlineHits.items[0].font.underline = \"Single\";\nboldHits.items[0].font.bold = true;\nawait context.sync();\n\n// Paragraph 3: currently holds only the (empty/collapsed) \"_GoBack\" bookmark.\n// Rebuild it as: plain lead-in text, then an italic \"Dedicated verification\n// permission\" run wrapped in the \"_GoBack\" bookmark, then plain trailing text.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst targetParagraph = paragraphs.items[2];\n\nconst middleRange = targetParagraph.insertText(\"Dedicated verification permission\", \"Start\");\nawait context.sync();\n\nmiddleRange.insertText(\n  \"Imap lotus decide ctr forward filter triggered id from subscribe started, users mailing group bayesian text list reply-to impression. \",\n  \"Before\"\n);\nmiddleRange.insertText(\n  \" challenge folder confirmation, mso into bring file users, lines font-size thank-you bug.\",\n  \"After\"\n);\nawait context.sync();\n\nmiddleRange.font.italic = true;\nawait context.sync();\n\nmiddleRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Paragraph 2: split the single run into 5 runs, underlining\n# \"mailing transactional\" and bolding \"Column enhanced fixed\".\n$d = $word.ActiveDocument\n\n$r1 = $d.Content\n$null = $r1.Find.Execute(\"mailing transactional\")\n$r1.Font.Underline = 1\n\n$r2 = $d.Content\n$null = $r2.Find.Execute(\"Column enhanced fixed\")\n$r2.Bold = 1\n\n# Paragraph 3: currently holds only the (empty/collapsed) \"_GoBack\" bookmark.\n# Rebuild it as: plain lead-in text, then an italic \"Dedicated verification\n# permission\" run wrapped in the \"_GoBack\" bookmark, then plain trailing text.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n$p3 = $d.Paragraphs.Item(3)\n$middleText = \"Dedicated verification permission\"\n$p3.Range.InsertAfter($middleText)\n\n$midStart = $p3.Range.Start\n$midEnd = $midStart + $middleText.Length\n$mid = $d.Range($midStart, $midEnd)\n$mid.Font.Italic = 1\n$d.Bookmarks.Add(\"_GoBack\", $mid)\n\n$mid.InsertBefore(\"Imap lotus decide ctr forward filter triggered id from subscribe started, users mailing group bayesian text list reply-to impression. \")\n$mid.InsertAfter(\" challenge folder confirmation, mso into bring file users, lines font-size thank-you bug.\")\n"}
